$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that might look numeric (e.g. "1.004", "22.455.33")
# as plain text, matching the source inline-string cells exactly. A leading
# apostrophe forces Excel's text parser; resetting Style to 'Normal' afterwards
# drops the quotePrefix style Excel would otherwise stamp on the cell, so the
# cell ends up with no explicit style (s attribute), same as the target cells.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '22.455.33'
$ws.Range('E2').Value = '  +9.00%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.604.66'
$ws.Range('E3').Value = '  +8.30%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.004'
$ws.Range('E4').Value = '  -0.70%  '

# Row 5
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range('D5') '305.71'
$ws.Range('E5').Value = '  +8.93%  '

# Row 6
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Range('D6') '0.9914'
$ws.Range('E6').Value = '  +2.07%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.3696'
$ws.Range('E7').Value = '  +0.88%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.3405'
$ws.Range('E8').Value = '  +10.02%  '

# Row 9
Set-TextValue $ws.Range('D9') '42.40'
$ws.Range('E9').Value = '  +5.52%  '

# Row 10
Set-TextValue $ws.Range('D10') '1.142'
$ws.Range('E10').Value = '  +7.24%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.07082'
$ws.Range('E11').Value = '  +5.99%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.9998'
$ws.Range('E12').Value = '  -0.63%  '

# Row 13
Set-TextValue $ws.Range('D13') '19.81'
$ws.Range('E13').Value = '  +8.86%  '

# Row 14
Set-TextValue $ws.Range('D14') '5.958'
$ws.Range('E14').Value = '  +7.49%  '

# Row 15
Set-TextValue $ws.Range('D15') '6.658'
$ws.Range('E15').Value = '  +6.89%  '

# Row 16
Set-TextValue $ws.Range('D16') '0.00001094'
$ws.Range('E16').Value = '  +5.96%  '

# Row 17
Set-TextValue $ws.Range('D17') '1.601.23'
$ws.Range('E17').Value = '  +7.89%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.9917'
$ws.Range('E18').Value = '  +2.12%  '

# Row 19
Set-TextValue $ws.Range('D19') '0.06823'
$ws.Range('E19').Value = '  +14.53%  '

# Row 20
Set-TextValue $ws.Range('D20') '78.16'
$ws.Range('E20').Value = '  +11.70%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D21') '6.054'
$ws.Range('E21').Value = '  +9.66%  '

# Row 22
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D22') '16.17'
$ws.Range('E22').Value = '  +10.95%  '

# Row 23
Set-TextValue $ws.Range('D23') '11.88'
$ws.Range('E23').Value = '  +7.03%  '

# Row 24
Set-TextValue $ws.Range('D24') '22.438.20'
$ws.Range('E24').Value = '  +8.70%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.410'
$ws.Range('E25').Value = '  +6.01%  '

# Row 26
Set-TextValue $ws.Range('D26') '2.550'
$ws.Range('E26').Value = '  +19.58%  '

# Row 27
Set-TextValue $ws.Range('D27') '151.30'
$ws.Range('E27').Value = '  +6.21%  '

# Row 28
Set-TextValue $ws.Range('D28') '19.63'
$ws.Range('E28').Value = '  +13.13%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.783.10'
$ws.Range('E29').Value = '  +8.34%  '

# Row 30
Set-TextValue $ws.Range('D30') '121.03'
$ws.Range('E30').Value = '  +5.69%  '

# Row 31
Set-TextValue $ws.Range('D31') '4.196'
$ws.Range('E31').Value = '  +6.18%  '

# Row 32
Set-TextValue $ws.Range('D32') '6.165'
$ws.Range('E32').Value = '  +21.99%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.9566'
$ws.Range('E33').Value = '  +15.00%  '

# Row 34
Set-TextValue $ws.Range('D34') '0.08308'
$ws.Range('E34').Value = '  +3.43%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.639'
$ws.Range('E35').Value = '  +6.06%  '

# Row 36
Set-TextValue $ws.Range('D36') '5.305'
$ws.Range('E36').Value = '  +11.10%  '

# Row 37
Set-TextValue $ws.Range('D37') '12.02'
$ws.Range('E37').Value = '  +14.70%  '

# Row 38
Set-TextValue $ws.Range('D38') '1.268'
$ws.Range('E38').Value = '  +4.11%  '

# Row 39
Set-TextValue $ws.Range('D39') '8.645'
$ws.Range('E39').Value = '  +12.13%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.06140'
$ws.Range('E40').Value = '  +5.62%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.02237'
$ws.Range('E41').Value = '  +9.03%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.2034'
$ws.Range('E42').Value = '  +7.70%  '

# Row 43
Set-TextValue $ws.Range('D43') '0.9916'
$ws.Range('E43').Value = '  +2.13%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.5941'
$ws.Range('E44').Value = '  +11.57%  '

# Row 45
Set-TextValue $ws.Range('D45') '3.852'
$ws.Range('E45').Value = '  +8.64%  '

# Row 46
Set-TextValue $ws.Range('D46') '13.13'
$ws.Range('E46').Value = '  +6.59%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.5719'
$ws.Range('E47').Value = '  +9.41%  '

# Row 48
Set-TextValue $ws.Range('D48') '127.86'
$ws.Range('E48').Value = '  +7.20%  '

# Row 49
Set-TextValue $ws.Range('D49') '1.988'
$ws.Range('E49').Value = '  +8.38%  '

# Row 50
Set-TextValue $ws.Range('D50') '0.06816'
$ws.Range('E50').Value = '  +4.57%  '

# Row 51
Set-TextValue $ws.Range('D51') '74.12'
$ws.Range('E51').Value = '  +8.95%  '
